$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(800, 50, 0.6, 0.7, 3432, 0),
    @(800, 50, 0.6, 0.7, 2551, 0),
    @(800, 50, 0.6, 0.7, 3411, 0),
    @(800, 50, 0.6, 0.7, 3997, 0)
)

$startRow = 301
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowValues[$c]
    }
}
